$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "2" rows (demand2, net2, pv2, bat2, CHP2) - these are currently
# rows 3, 5, 7, 9, 11. Delete bottom-up so row indices above stay valid.
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(3).Delete()
